$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 05:03"

# Update Pakistan row (row 26)
$ws.Range("B26").Value = 22550
$ws.Range("C26").Value = 501
$ws.Range("D26").Value = 6217
$ws.Range("E26").Value = 15807
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 526

# Update Venezuela row (row 127)
$ws.Range("B127").Value = 367
$ws.Range("C127").Value = 6
$ws.Range("D127").Value = 164

# Swap Curazao (row 198) and Dominica (row 199)
$ws.Range("A198").Value = "Dominica"
$ws.Range("B198").Value = 16
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 2
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("B199").Value = 16
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 13
$ws.Range("E199").Value = 2
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1

# Swap Seychelles (row 205) and Montserrat (row 206)
$ws.Range("A205").Value = "Montserrat"
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 7
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

$ws.Range("A206").Value = "Seychelles"
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 0
